$wb = $excel.ActiveWorkbook

# --- Sheet: Resource Utilization ---
$wsRU = $wb.Worksheets.Item("Resource Utilization")
$wsRU.Range("B2").Value = 100
$wsRU.Range("B3").Value = 45.35

# --- Sheet: Activity Times ---
$wsAT = $wb.Worksheets.Item("Activity Times")

# Row 2
$wsAT.Cells.Item(2, 1).Value = "Review AM using Asset Change Tracker (5.5.13.1)"
$wsAT.Cells.Item(2, 2).Value = "Activity Step"
$wsAT.Cells.Item(2, 3).Value = 146
$wsAT.Cells.Item(2, 4).Value = 146
$wsAT.Cells.Item(2, 5).Value = 6
$wsAT.Cells.Item(2, 6).Value = 13
$wsAT.Cells.Item(2, 7).Value = 10.42
$wsAT.Cells.Item(2, 8).Value = 0
$wsAT.Cells.Item(2, 9).Value = 0
$wsAT.Cells.Item(2, 10).Value = 0
$wsAT.Cells.Item(2, 11).Value = 0

# Row 3
$wsAT.Cells.Item(3, 1).Value = "Complete /Accurate?"
$wsAT.Cells.Item(3, 2).Value = "Gateway"
$wsAT.Cells.Item(3, 3).Value = 146
$wsAT.Cells.Item(3, 4).Value = 146
$wsAT.Cells.Item(3, 5).Value = 1
$wsAT.Cells.Item(3, 6).Value = 1
$wsAT.Cells.Item(3, 7).Value = 1
$wsAT.Cells.Item(3, 8).Value = 0
$wsAT.Cells.Item(3, 9).Value = 0
$wsAT.Cells.Item(3, 10).Value = 0
$wsAT.Cells.Item(3, 11).Value = 0

# Row 4
$wsAT.Cells.Item(4, 1).Value = "Work with REO RPO to Correct (5.5.13.3)"
$wsAT.Cells.Item(4, 2).Value = "Activity Step"
$wsAT.Cells.Item(4, 3).Value = 60
$wsAT.Cells.Item(4, 4).Value = 60
$wsAT.Cells.Item(4, 5).Value = 24
$wsAT.Cells.Item(4, 6).Value = 165
$wsAT.Cells.Item(4, 7).Value = 117.12
$wsAT.Cells.Item(4, 8).Value = 0
$wsAT.Cells.Item(4, 9).Value = 0
$wsAT.Cells.Item(4, 10).Value = 0
$wsAT.Cells.Item(4, 11).Value = 0

# Row 5
$wsAT.Cells.Item(5, 1).Value = "Note Accuracy in Asset Change Tracker (5.5.13.2)"
$wsAT.Cells.Item(5, 2).Value = "Activity Step"
$wsAT.Cells.Item(5, 3).Value = 147
$wsAT.Cells.Item(5, 4).Value = 147
$wsAT.Cells.Item(5, 5).Value = 1
$wsAT.Cells.Item(5, 6).Value = 6
$wsAT.Cells.Item(5, 7).Value = 3.84
$wsAT.Cells.Item(5, 8).Value = 1698
$wsAT.Cells.Item(5, 9).Value = 0
$wsAT.Cells.Item(5, 10).Value = 165
$wsAT.Cells.Item(5, 11).Value = 22.64

# Row 6
$wsAT.Cells.Item(6, 1).Value = "Create/Post Journal Entries (5.5.13.4)"
$wsAT.Cells.Item(6, 2).Value = "Stop"
$wsAT.Cells.Item(6, 3).Value = 144
$wsAT.Cells.Item(6, 4).Value = 144
$wsAT.Cells.Item(6, 5).Value = 3
$wsAT.Cells.Item(6, 6).Value = 6
$wsAT.Cells.Item(6, 7).Value = 5.15
$wsAT.Cells.Item(6, 8).Value = 0
$wsAT.Cells.Item(6, 9).Value = 0
$wsAT.Cells.Item(6, 10).Value = 0
$wsAT.Cells.Item(6, 11).Value = 0

# Row 7
$wsAT.Cells.Item(7, 1).Value = "Stop"
$wsAT.Cells.Item(7, 2).Value = "Unknown"
$wsAT.Cells.Item(7, 3).Value = 144
$wsAT.Cells.Item(7, 4).Value = 144
$wsAT.Cells.Item(7, 5).Value = 0
$wsAT.Cells.Item(7, 6).Value = 0
$wsAT.Cells.Item(7, 7).Value = 0
$wsAT.Cells.Item(7, 8).Value = 0
$wsAT.Cells.Item(7, 9).Value = 0
$wsAT.Cells.Item(7, 10).Value = 0
$wsAT.Cells.Item(7, 11).Value = 0
